# Helper: force a cell (or range) to hold a literal TEXT value, even if the
# value looks numeric (e.g. "510810" or "68.43"), without leaving any
# leftover number-format style on the cell.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# The workbook currently has sheets:
#   2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 总计
# Target layout:
#   2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 2022-Q1, 总计
#
# The existing "总计" sheet becomes the new "2022-Q1" sheet (its data is
# replaced with the 2022-Q1 fund holding table), and a fresh copy of the
# original "总计" sheet is appended at the end, keeping the "总计" name and
# receiving an extra row for the new 2022-Q1 summary data.
# ---------------------------------------------------------------------

$oldTotal = $wb.Worksheets.Item("总计")

# Duplicate the current "总计" sheet right after itself; this copy will
# become the new "总计" summary sheet.
$oldTotal.Copy($null, $oldTotal)
$newTotal = $wb.Worksheets.Item($oldTotal.Index + 1)

# Rename sheets into their final positions/names.
$oldTotal.Name = "2022-Q1"
$newTotal.Name = "总计"

# =======================================================================
# 1) Fill the "2022-Q1" sheet with the fund holding table.
# =======================================================================
$ws = $oldTotal

# Extend the header row with the new columns (E1:H1), copying the
# existing header style (s="2") from D1.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1:H1").PasteSpecial(-4122) | Out-Null

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Copy the styled column-A cell (s="2") down to the two brand new rows
# (rows 8 and 9) before writing data into them.
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A7:A9").PasteSpecial(-4122) | Out-Null

$fundRows = @(
    @{ Row = 2; Idx = 0; Code = "510810"; Name = "汇添富中证上海国企ETF";     Size = "68.43"; Pos = "99.71"; Pct = "3.77";  Value = "2.5798"; Rank = 8 },
    @{ Row = 3; Idx = 1; Code = "161721"; Name = "招商沪深300地产等权重指数"; Size = "9.97";  Pos = "94.51"; Pct = "10.26"; Value = "1.0229"; Rank = 7 },
    @{ Row = 4; Idx = 2; Code = "512200"; Name = "南方中证全指房地产ETF";     Size = "28.63"; Pos = "99.85"; Pct = "1.74";  Value = "0.4982"; Rank = 9 },
    @{ Row = 5; Idx = 3; Code = "160218"; Name = "国泰国证房地产行业指数";     Size = "6.35";  Pos = "95.04"; Pct = "1.92";  Value = "0.1219"; Rank = 8 },
    @{ Row = 6; Idx = 4; Code = "002585"; Name = "建信兴利灵活配置混合";       Size = "2.04";  Pos = "61.22"; Pct = "5.71";  Value = "0.1165"; Rank = 3 },
    @{ Row = 7; Idx = 5; Code = "003831"; Name = "建信鑫瑞回报灵活配置混合";   Size = "1.80";  Pos = "70.33"; Pct = "6.17";  Value = "0.1111"; Rank = 1 },
    @{ Row = 8; Idx = 6; Code = "160628"; Name = "鹏华中证800地产指数（LOF）"; Size = "3.38";  Pos = "94.35"; Pct = "2.46";  Value = "0.0831"; Rank = 9 },
    @{ Row = 9; Idx = 7; Code = "515060"; Name = "华夏中证全指房地产ETF";     Size = "2.37";  Pos = "98.82"; Pct = "1.72";  Value = "0.0408"; Rank = 9 }
)

foreach ($r in $fundRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Idx
    Set-TextValue $ws.Range("B$row") $r.Code
    $ws.Range("C$row").Value = $r.Name
    Set-TextValue $ws.Range("D$row") $r.Size
    Set-TextValue $ws.Range("E$row") $r.Pos
    Set-TextValue $ws.Range("F$row") $r.Pct
    Set-TextValue $ws.Range("G$row") $r.Value
    $ws.Range("H$row").Value = $r.Rank
}

# =======================================================================
# 2) Fill the new "总计" sheet with the updated date/count/value summary,
#    inserting the new 2022-Q1 row at the top.
# =======================================================================
$ws2 = $newTotal

# Copy the styled column-A cell (s="2") down to the newly needed row 7.
$ws2.Range("A6").Copy() | Out-Null
$ws2.Range("A7").PasteSpecial(-4122) | Out-Null

$summaryRows = @(
    @{ Row = 2; Idx = 0; Date = "2022-Q1"; Count = 8; Value = 4.57 },
    @{ Row = 3; Idx = 1; Date = "2021-Q4"; Count = 8; Value = 7.22 },
    @{ Row = 4; Idx = 2; Date = "2021-Q3"; Count = 5; Value = 4.16 },
    @{ Row = 5; Idx = 3; Date = "2021-Q2"; Count = 6; Value = 4.58 },
    @{ Row = 6; Idx = 4; Date = "2021-Q1"; Count = 2; Value = 3.72 },
    @{ Row = 7; Idx = 5; Date = "2020-Q4"; Count = 2; Value = 4.29 }
)

foreach ($r in $summaryRows) {
    $row = $r.Row
    $ws2.Range("A$row").Value = $r.Idx
    $ws2.Range("B$row").Value = $r.Date
    $ws2.Range("C$row").Value = $r.Count
    $ws2.Range("D$row").Value = $r.Value
}

# Restore the original active sheet (2020-Q4) so the workbook's selected
# tab matches the pre-edit state instead of defaulting to the last-touched
# sheet.
$wb.Worksheets.Item(1).Activate()
